$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the example values in column A (rows 2-6) from "exampleN" to "eN"
$ws.Range("A2").Value = "e1"
$ws.Range("A3").Value = "e2"
$ws.Range("A4").Value = "e3"
$ws.Range("A5").Value = "e4"
$ws.Range("A6").Value = "e5"

# Update the active selection from A4 to B2
$ws.Range("B2").Select()
